# Automatische test-sync: 2025-08-28 20:37:50
# Append the new "Opvolging retour" log entry to the Logs sheet and
# bump the matching Dashboard summary count.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")

# --- Append new row 16 to the Logs sheet ---
$wsLogs.Range("A16").Value = "Opvolging retour"
$wsLogs.Range("B16").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("D16").Value = "Retour / Terugbetaling"
$wsLogs.Range("F16").Value = "2025-08-28 20:37:22"
$wsLogs.Range("G16").Value = "Nee"
$wsLogs.Range("H16").Value = "Ja"
$wsLogs.Range("I16").Value = "Nee"
$wsLogs.Range("J16").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$dConditions = $wsLogs.Range("D2:D15").FormatConditions
for ($i = 1; $i -le $dConditions.Count; $i++) {
    $dConditions.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D16"))
}

$gConditions = $wsLogs.Range("G2:G15").FormatConditions
for ($i = 1; $i -le $gConditions.Count; $i++) {
    $gConditions.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G16"))
}

$hConditions = $wsLogs.Range("H2:H15").FormatConditions
for ($i = 1; $i -le $hConditions.Count; $i++) {
    $hConditions.Item($i).ModifyAppliesToRange($wsLogs.Range("H2:H16"))
}

$iConditions = $wsLogs.Range("I2:I15").FormatConditions
for ($i = 1; $i -le $iConditions.Count; $i++) {
    $iConditions.Item($i).ModifyAppliesToRange($wsLogs.Range("I2:I16"))
}

$jConditions = $wsLogs.Range("J2:J15").FormatConditions
for ($i = 1; $i -le $jConditions.Count; $i++) {
    $jConditions.Item($i).ModifyAppliesToRange($wsLogs.Range("J2:J16"))
}

# --- Update the Dashboard summary count for "Retour / Terugbetaling" ---
$wsDashboard = $wb.Worksheets.Item("Dashboard")
$wsDashboard.Range("B2").Value = 15
